# Fixing issues in segment\final.m
#
# Add three new variable rows to the bottom of the "Variables" table on
# Sheet1 (rows 41-43), following the existing layout: column B holds the
# (longer) description, column C holds the short variable name/token.
#
#   row 41: k_H      = Stiffness of pod shell (N/m^1.5)
#   row 42: T        = Length of time segments (s)
#   row 43: B_m_min  = Minimum value of perpendicular vector (m/s)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the three new rows by copying the format of the last existing data
# row (40) downward, so the new B/C cells pick up the same cell styles
# (s="2" / s="3") already used throughout the table, instead of Excel's
# plain default formatting.
$ws.Rows(40).Copy()
$null = $ws.Rows("41:43").Insert(-4121, -4163)   # xlShiftDown, xlFormatFromLeftOrAbove
$excel.CutCopyMode = $false

# Fill in the name (column C) before the description (column B) on each
# row, matching the original authoring order of the shared-string table.
$ws.Cells.Item(41, 3).Value = "k_H = "
$ws.Cells.Item(41, 2).Value = "Stiffness of pod shell (N/m^1.5)"

$ws.Cells.Item(42, 3).Value = "T = "
$ws.Cells.Item(42, 2).Value = "Length of time segments (s)"

$ws.Cells.Item(43, 3).Value = "B_m_min = "
$ws.Cells.Item(43, 2).Value = "Minimum value of perpendicular vector (m/s)"

# Leave the selection where the author ended up after typing the new rows.
$null = $ws.Range("B44").Select()
